$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StepperPage")

# Update the login cells for the stepper-page automation test data
$ws.Range("B2").Value = "automationtest@pixentia.com"
$ws.Range("D2").Value = "automation"

# Turn the username cell into a mailto hyperlink
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:automationtest@pixentia.com")
$ws.Range("B2").Style = "Hyperlink"

# Move the active selection to D2
$ws.Range("D2").Select() | Out-Null
